$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M12").Value = 12246.22

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F12").Value = 12246.22
$wsMensual.Range("F24").Value = 23459.63
$wsMensual.Columns.Item(6).ColumnWidth = 13.1666666666667

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D16").Value = 23459.63
$wsCumplimiento.Range("E16").Value = 15296.91
$wsCumplimiento.Range("F16").Value = 0.6053076461417867
$wsCumplimiento.Range("D19").Value = 23459.63
$wsCumplimiento.Range("E19").Value = 34763.37386304603
$wsCumplimiento.Range("F19").Value = 0.4029271669868232
$wsCumplimiento.Columns.Item(6).ColumnWidth = 23.1666666666667
